# Updated cryptos list with latest prices and 1h volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.599.06'
$ws.Range("E2").Value = '  +3.46%  '
$ws.Range("D3").Value = '''2.648.33'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = '''569.58'
$ws.Range("E5").Value = '  +6.45%  '
$ws.Range("D6").Value = '''147.08'
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '''0.611'
$ws.Range("E8").Value = '  +6.95%  '
$ws.Range("D9").Value = '''2.671.68'
$ws.Range("E9").Value = '  +1.62%  '
$ws.Range("D10").Value = '''6.84'
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  +4.86%  '
$ws.Range("E12").Value = '  +6.30%  '
$ws.Range("D13").Value = '''0.344'
$ws.Range("E13").Value = '  +3.20%  '
$ws.Range("D14").Value = '''3.116.30'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").Value = '''60.586.65'
$ws.Range("E15").Value = '  +3.58%  '
$ws.Range("E16").Value = '  +5.49%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '''2.670.79'
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '''0.0000138'
$ws.Range("E18").Value = '  +4.45%  '
$ws.Range("E19").Value = '  +3.53%  '
$ws.Range("D20").Value = '''344.64'
$ws.Range("E20").Value = '  +2.89%  '
$ws.Range("D21").Value = '''10.49'
$ws.Range("E21").Value = '  +3.30%  '
$ws.Range("E22").Value = '  +2.49%  '
$ws.Range("D23").Value = '''5.84'
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").Value = '''66.64'
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").Value = '  +5.71%  '
$ws.Range("E27").Value = '  +1.77%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").Value = '''7.42'
$ws.Range("E29").Value = '  +4.21%  '
$ws.Range("D30").Value = '''0.0₃0787'
$ws.Range("E30").Value = '  +6.79%  '
$ws.Range("D31").Value = '''0.998'
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").Value = '''6.33'
$ws.Range("E32").Value = '  +7.61%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''1.72'
$ws.Range("E33").Value = '  +4.89%  '
$ws.Range("D34").Value = '''19.25'
$ws.Range("E34").Value = '  +1.75%  '
$ws.Range("D35").Value = '''154.73'
$ws.Range("E35").Value = '  +2.37%  '
$ws.Range("E36").Value = '  +5.18%  '
$ws.Range("E37").Value = '  +7.77%  '
$ws.Range("D38").Value = '''0.908'
$ws.Range("E38").Value = '  +6.43%  '
$ws.Range("D39").Value = '''0.911'
$ws.Range("E39").Value = '  +11.99%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("E41").Value = '  +7.32%  '
$ws.Range("D42").Value = '''304.08'
$ws.Range("E42").Value = '  +7.80%  '
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("D44").Value = '''0.993'
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '''0.0984'
$ws.Range("E45").Value = '  +4.89%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.607'
$ws.Range("E46").Value = '  +2.22%  '
$ws.Range("E47").Value = '  +4.51%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''19.54'
$ws.Range("E48").Value = '  +2.56%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '''128.89'
$ws.Range("E49").Value = '  +12.92%  '
$ws.Range("D50").Value = '''10.67'
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = '''0.0236'
$ws.Range("E51").Value = '  +5.06%  '
